$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Bùi Anh Dũng"
$ws.Cells.Item(4, 3).Value = 12
$ws.Cells.Item(4, 4).Value = 2025
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 10000000
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
